# Auto-generated Excel COM-interop script
# Applies weekly crime-data refresh edit (Volume 31 Number 43 -> 44,
# report week 10/21/2024-10/27/2024 -> 10/28/2024-11/3/2024, and updated crime stats).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number + report week dates) ---
$ws.Range("C1").Value = "Volume 31   Number  44"
$ws.Range("C6").Value = "Report Covering the Week  10/28/2024  Through  11/3/2024"

# --- Column width adjustments (columns H, I, J narrower) ---
$ws.Columns.Item(8).ColumnWidth = 6.168446
$ws.Columns.Item(9).ColumnWidth = 6.168446
$ws.Columns.Item(10).ColumnWidth = 6.168446

# --- Donor cells (row 33) used to clone formats for N/A <-> numeric toggles ---
# Row 33 formatting is untouched by this edit, so it is safe to use as a format source.
$donorText0   = $ws.Range("C33")   # style 13, shared string "0"
$donorTextNA  = $ws.Range("E33")   # style 13, shared string "***.*"
$donorNumber  = $ws.Range("I33")   # style 14, plain whole number
$donorDecimal = $ws.Range("K33")   # style 15, decimal percentage

# --- Cells that flip from numeric to the blank-marker text style (copy format+value) ---
$donorText0.Copy($ws.Range("F14"))
$donorText0.Copy($ws.Range("G27"))
$donorTextNA.Copy($ws.Range("H27"))
$donorText0.Copy($ws.Range("C28"))
$donorText0.Copy($ws.Range("D31"))
$donorTextNA.Copy($ws.Range("E31"))

# --- Cells that flip from the blank-marker text style back to numeric (copy format, then set value) ---
$donorNumber.Copy($ws.Range("C29"))
$ws.Range("C29").Value = 2
$donorNumber.Copy($ws.Range("D29"))
$ws.Range("D29").Value = 2
$donorDecimal.Copy($ws.Range("E29"))
$ws.Range("E29").Value = 0
$donorNumber.Copy($ws.Range("G29"))
$ws.Range("G29").Value = 2
$donorDecimal.Copy($ws.Range("H29"))
$ws.Range("H29").Value = 0
$donorNumber.Copy($ws.Range("C30"))
$ws.Range("C30").Value = 1
$donorNumber.Copy($ws.Range("D30"))
$ws.Range("D30").Value = 1
$donorDecimal.Copy($ws.Range("E30"))
$ws.Range("E30").Value = 0
$donorNumber.Copy($ws.Range("G30"))
$ws.Range("G30").Value = 1
$donorDecimal.Copy($ws.Range("H30"))
$ws.Range("H30").Value = 0

# --- Plain value updates (style unchanged) ---
$ws.Range("I15").Value = 20
$ws.Range("K15").Value = 53.846153846153
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 66.666666666666
$ws.Range("N15").Value = -25.925925925925
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -75
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -12.5
$ws.Range("I16").Value = 101
$ws.Range("J16").Value = 124
$ws.Range("K16").Value = -18.548387096774
$ws.Range("L16").Value = -36.075949367088
$ws.Range("M16").Value = -54.090909090909
$ws.Range("N16").Value = -88.117647058823
$ws.Range("C17").Value = 8
$ws.Range("E17").Value = 33.333333333333
$ws.Range("F17").Value = 22
$ws.Range("G17").Value = 25
$ws.Range("H17").Value = -12
$ws.Range("I17").Value = 246
$ws.Range("J17").Value = 201
$ws.Range("K17").Value = 22.388059701492
$ws.Range("L17").Value = 47.305389221556
$ws.Range("M17").Value = 192.857142857143
$ws.Range("N17").Value = 17.142857142857
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -80
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = -10.526315789473
$ws.Range("I18").Value = 184
$ws.Range("J18").Value = 224
$ws.Range("K18").Value = -17.857142857142
$ws.Range("L18").Value = -7.070707070707
$ws.Range("M18").Value = -22.362869198312
$ws.Range("N18").Value = -84.979591836734
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = -60
$ws.Range("F19").Value = 36
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = -25
$ws.Range("I19").Value = 458
$ws.Range("J19").Value = 520
$ws.Range("K19").Value = -11.923076923076
$ws.Range("L19").Value = -19.366197183098
$ws.Range("M19").Value = 24.119241192411
$ws.Range("N19").Value = -11.583011583011
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = -30
$ws.Range("F20").Value = 32
$ws.Range("G20").Value = 22
$ws.Range("H20").Value = 45.454545454545
$ws.Range("I20").Value = 326
$ws.Range("J20").Value = 302
$ws.Range("K20").Value = 7.947019867549
$ws.Range("L20").Value = 42.982456140350
$ws.Range("M20").Value = 65.482233502538
$ws.Range("N20").Value = -92.174747959673
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 40
$ws.Range("E21").Value = -40
$ws.Range("F21").Value = 124
$ws.Range("G21").Value = 130
$ws.Range("H21").Value = -4.615384615384
$ws.Range("I21").Value = 1339
$ws.Range("J21").Value = 1385
$ws.Range("K21").Value = -3.321299638989
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 19.234194122885
$ws.Range("N21").Value = -80.876892316481
$ws.Range("C22").Value = 4
$ws.Range("E22").Value = 300
$ws.Range("F22").Value = 6
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 34
$ws.Range("J22").Value = 36
$ws.Range("K22").Value = -5.555555555555
$ws.Range("L22").Value = 36
$ws.Range("M22").Value = 61.904761904761
$ws.Range("D23").Value = 2
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = -83.333333333333
$ws.Range("J23").Value = 67
$ws.Range("K23").Value = -28.358208955223
$ws.Range("L23").Value = 6.666666666666
$ws.Range("M23").Value = 54.838709677419
$ws.Range("C24").Value = 26
$ws.Range("E24").Value = -3.703703703703
$ws.Range("F24").Value = 97
$ws.Range("H24").Value = -3.960396039603
$ws.Range("I24").Value = 1016
$ws.Range("J24").Value = 1207
$ws.Range("K24").Value = -15.824357912179
$ws.Range("L24").Value = -18.914604948124
$ws.Range("M24").Value = 19.388954171562
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 43
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = 22.857142857142
$ws.Range("I25").Value = 434
$ws.Range("J25").Value = 443
$ws.Range("K25").Value = -2.031602708803
$ws.Range("L25").Value = 15.425531914893
$ws.Range("C26").Value = 7
$ws.Range("E26").Value = -30
$ws.Range("F26").Value = 40
$ws.Range("G26").Value = 43
$ws.Range("H26").Value = -6.976744186046
$ws.Range("I26").Value = 446
$ws.Range("J26").Value = 409
$ws.Range("K26").Value = 9.046454767726
$ws.Range("L26").Value = 17.060367454068
$ws.Range("M26").Value = 28.901734104046
$ws.Range("F27").Value = 3
$ws.Range("I27").Value = 31
$ws.Range("K27").Value = 34.782608695652
$ws.Range("L27").Value = 6.896551724137
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = -57.142857142857
$ws.Range("J28").Value = 42
$ws.Range("K28").Value = -19.047619047619
$ws.Range("L28").Value = -33.333333333333
$ws.Range("F29").Value = 2
$ws.Range("I29").Value = 6
$ws.Range("J29").Value = 4
$ws.Range("K29").Value = 50
$ws.Range("L29").Value = -33.333333333333
$ws.Range("M29").Value = 20
$ws.Range("N29").Value = -40
$ws.Range("I30").Value = 5
$ws.Range("J30").Value = 3
$ws.Range("K30").Value = 66.666666666666
$ws.Range("L30").Value = -28.571428571428
$ws.Range("M30").Value = 0
$ws.Range("N30").Value = -50
